$d = $word.ActiveDocument

# --- #2: array average explanation paragraph ---
$old2 = "In this recursive algorithm, we pass along A, the array of integers, and n, the position in the array we are. In the base case we return the first element in the list, otherwise we calculate the average by recursively calling our method with the previous element in the list (which gives us the sum up until the previous element), multiplied by that element, plus the current element, and divide by n to average that many values. At our base case when n==1, we return the first element, then move up the stack until we are able to make our final calculation, the average."
$new2 = "In this recursive algorithm, we pass along A, the array of integers, and n, the position in the array we are, initially set to the length of the list. We recursively call our method with the previous element in the list. When we reach our base case when n==1, we return the first element of the list, and get the sum of the elements from the beginning of the array to n-2, we then add the element at n-1 and divide by n to get the average of all the elements."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# --- #2: pseudocode method signature comment ---
$old3 = "Method arrayAvg(A, n)"
$new3 = "Method arrayAvg(A, n) //where A is the array and n is the length of the array A"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# --- #2: pseudocode return line, A[n] -> A[n-1] ---
$old4 = "return (arrayAvg(A, n-1)*(n-1) + A[n])/n"
$new4 = "return (arrayAvg(A, n-1)*(n-1) + A[n-1])/n"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)

# --- #3: binary search explanation paragraph wording tweak ---
$old5 = "Each time we recursively call binary search we halve the search space, so at worst the algorithm could halve the search space until it is at one element left. This depends on the number of elements n, and would occur within log(n) +1 calls, log(n) to search the space recursively and divide it repeatedly, and +1 to account for the initial call. "
$new5 = "Each time we recursively call binary search we halve the search space, so at worst the algorithm could halve the search space until there is only one element left. This depends on the number of elements n, and would occur within log(n) +1 calls, log(n) to search the space recursively by dividing it repeatedly, and +1 to account for the initial call. "
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)

# --- #5: insert a blank paragraph before the "5." heading ---
# Scope the search to after the gcd pseudocode so we don't match the "5." inside
# "EN.605.202.81" earlier in the document.
$anchor5 = $d.Content
$anchor5.Find.Execute("gcd(y, x%y)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng5 = $d.Range($anchor5.End, $d.Content.End)
$rng5.Find.Execute("5.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng5.InsertParagraphBefore()

# --- #5: fibonacci explanation paragraph, append additional sentence ---
$old6 = "Recursive fibonacci is defined with two base cases for values 0 and 1, and a recursive case when n>1."
$new6 = "Recursive fibonacci is defined with two base cases for values 0 and 1, and a recursive case when n>1. The recursive case involves the addition of the two previous fib elements, giving us the fibonacci sequence. Gfib involves passing starting parameters f0 and f1."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
